$d = $word.ActiveDocument

# Replace each occurrence of the word "most" (immediately followed by " paid")
# with "highest", while forcing Word to keep the replaced word in its own
# run (matching the author's edit, which shows each paragraph split into
# three runs, e.g. "...the ", "highest", " paid...").
#
# A plain Range.Text assignment (or Find.Execute replace) on this engine
# re-coalesces all same-formatted runs in a paragraph into a single run,
# so we briefly flip a character formatting property on the replacement
# range and then clear it again. That forces the run boundary to persist
# without leaving any visible formatting behind.

$searchFrom = 0
for ($i = 0; $i -lt 2; $i++) {
    $content = $d.Content.Text
    $idx = $content.IndexOf("most paid", $searchFrom)
    if ($idx -lt 0) {
        break
    }

    $target = $d.Range($idx, $idx + 4)
    $target.Bold = 1
    $target.Text = "highest"

    $newRun = $d.Range($idx, $idx + 7)
    $newRun.Bold = 0

    $searchFrom = $idx + 7
}
